$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 109.5  # H11: 60 -> 109.5
$ws.Cells.Item(11, 9).Value = 109.5  # I11: 60 -> 109.5
$ws.Cells.Item(11, 11).Value = 109.5  # K11: 60 -> 109.5
$ws.Cells.Item(11, 13).Value = 30.5  # M11: 80 -> 30.5
$ws.Cells.Item(15, 8).Value = 2854.7966  # H15: 3004.1606 -> 2854.7966
$ws.Cells.Item(15, 9).Value = 2854.7966  # I15: 3004.1606 -> 2854.7966
$ws.Cells.Item(15, 11).Value = 8564.389800000001  # K15: 9012.481800000001 -> 8564.389800000001
$ws.Cells.Item(15, 13).Value = -8395.389800000001  # M15: -8843.481800000001 -> -8395.389800000001
$ws.Cells.Item(33, 8).Value = 156.15  # H33: 167.22223 -> 156.15
$ws.Cells.Item(33, 9).Value = 163.66667  # I33: 177.0625 -> 163.66667
$ws.Cells.Item(33, 11).Value = 163.66667  # K33: 177.0625 -> 163.66667
$ws.Cells.Item(33, 13).Value = 65.33332999999999  # M33: 51.9375 -> 65.33332999999999
$ws.Cells.Item(41, 8).Value = 664.6667  # H41: 282.44446 -> 664.6667
$ws.Cells.Item(41, 9).Value = 165.71428  # I41: 151.33333 -> 165.71428
$ws.Cells.Item(41, 10).Value = 1363.2  # J41: 544.6667 -> 1363.2
$ws.Cells.Item(41, 11).Value = 165.71428  # K41: 151.33333 -> 165.71428
$ws.Cells.Item(41, 12).Value = 1363.2  # L41: 544.6667 -> 1363.2
$ws.Cells.Item(41, 13).Value = 274.28572  # M41: 288.66667 -> 274.28572
$ws.Cells.Item(41, 14).Value = -2243.2  # N41: -1424.6667 -> -2243.2
$ws.Cells.Item(50, 8).Value = 196  # H50: 198 -> 196
$ws.Cells.Item(50, 10).Value = 196  # J50: 198 -> 196
$ws.Cells.Item(50, 12).Value = 588  # L50: 594 -> 588
$ws.Cells.Item(50, 14).Value = -1538  # N50: -1544 -> -1538
$ws.Cells.Item(80, 8).Value = 405.86667  # H80: 391.64706 -> 405.86667
$ws.Cells.Item(80, 9).Value = 410.35715  # I80: 428.76923 -> 410.35715
$ws.Cells.Item(80, 10).Value = 343  # J80: 271 -> 343
$ws.Cells.Item(80, 11).Value = 1231.07145  # K80: 1286.30769 -> 1231.07145
$ws.Cells.Item(80, 12).Value = 1029  # L80: 813 -> 1029
$ws.Cells.Item(80, 13).Value = -233.0714499999999  # M80: -288.3076900000001 -> -233.0714499999999
$ws.Cells.Item(80, 14).Value = -3025  # N80: -2809 -> -3025
$ws.Cells.Item(83, 8).Value = 405.86667  # H83: 391.64706 -> 405.86667
$ws.Cells.Item(83, 9).Value = 410.35715  # I83: 428.76923 -> 410.35715
$ws.Cells.Item(83, 10).Value = 343  # J83: 271 -> 343
$ws.Cells.Item(83, 11).Value = 3693.21435  # K83: 3858.92307 -> 3693.21435
$ws.Cells.Item(83, 12).Value = 3087  # L83: 2439 -> 3087
$ws.Cells.Item(83, 13).Value = 1298.78565  # M83: 1133.07693 -> 1298.78565
$ws.Cells.Item(83, 14).Value = -13071  # N83: -12423 -> -13071
$ws.Cells.Item(111, 8).Value = 1539.1111  # H111: 1861.1428 -> 1539.1111
$ws.Cells.Item(111, 9).Value = 1792.3334  # I111: 2382 -> 1792.3334
$ws.Cells.Item(111, 10).Value = 1032.6666  # J111: 1166.6666 -> 1032.6666
$ws.Cells.Item(111, 11).Value = 5377.0002  # K111: 7146 -> 5377.0002
$ws.Cells.Item(111, 12).Value = 3097.9998  # L111: 3499.9998 -> 3097.9998
$ws.Cells.Item(111, 13).Value = -2310.0002  # M111: -4079 -> -2310.0002
$ws.Cells.Item(111, 14).Value = -9231.9998  # N111: -9633.9998 -> -9231.9998
$ws.Cells.Item(112, 8).Value = 1936.2307  # H112: 1905.1875 -> 1936.2307
$ws.Cells.Item(112, 10).Value = 2006.4546  # J112: 1955.9286 -> 2006.4546
$ws.Cells.Item(112, 12).Value = 6019.3638  # L112: 5867.7858 -> 6019.3638
$ws.Cells.Item(112, 14).Value = -8235.363799999999  # N112: -8083.7858 -> -8235.363799999999
$ws.Cells.Item(115, 8).Value = 1749.2  # H115: 2534.6667 -> 1749.2
$ws.Cells.Item(115, 9).Value = 436.5  # I115: 427 -> 436.5
$ws.Cells.Item(115, 10).Value = 7000  # J115: 6750 -> 7000
$ws.Cells.Item(115, 11).Value = 1309.5  # K115: 1281 -> 1309.5
$ws.Cells.Item(115, 12).Value = 21000  # L115: 20250 -> 21000
$ws.Cells.Item(115, 13).Value = 257.5  # M115: 286 -> 257.5
$ws.Cells.Item(115, 14).Value = -24134  # N115: -23384 -> -24134
$ws.Cells.Item(129, 8).Value = 930.5  # H129: 939.8 -> 930.5
$ws.Cells.Item(129, 9).Value = 755.38464  # I129: 756.0769 -> 755.38464
$ws.Cells.Item(129, 10).Value = 1689.3334  # J129: 2134 -> 1689.3334
$ws.Cells.Item(129, 11).Value = 2266.15392  # K129: 2268.2307 -> 2266.15392
$ws.Cells.Item(129, 12).Value = 5068.0002  # L129: 6402 -> 5068.0002
$ws.Cells.Item(129, 13).Value = 2733.84608  # M129: 2731.7693 -> 2733.84608
$ws.Cells.Item(129, 14).Value = -15068.0002  # N129: -16402 -> -15068.0002
$ws.Cells.Item(137, 8).Value = 928798.9399999999  # H137: 1002944.8 -> 928798.9399999999
$ws.Cells.Item(137, 9).Value = 2501571.8  # I137: 2779407.2 -> 2501571.8
$ws.Cells.Item(137, 10).Value = 3638.5  # J137: 3684.6562 -> 3638.5
$ws.Cells.Item(137, 11).Value = 7504715.399999999  # K137: 8338221.600000001 -> 7504715.399999999
$ws.Cells.Item(137, 12).Value = 10915.5  # L137: 11053.9686 -> 10915.5
$ws.Cells.Item(137, 13).Value = -7502165.399999999  # M137: -8335671.600000001 -> -7502165.399999999
$ws.Cells.Item(137, 14).Value = -16015.5  # N137: -16153.9686 -> -16015.5
$ws.Cells.Item(138, 8).Value = 3351.3958  # H138: 3407.6 -> 3351.3958
$ws.Cells.Item(138, 9).Value = 3697.7  # I138: 3506.5454 -> 3697.7
$ws.Cells.Item(138, 10).Value = 3104.0356  # J138: 3286.6667 -> 3104.0356
$ws.Cells.Item(138, 11).Value = 11093.1  # K138: 10519.6362 -> 11093.1
$ws.Cells.Item(138, 12).Value = 9312.106800000001  # L138: 9860.000100000001 -> 9312.106800000001
$ws.Cells.Item(138, 13).Value = -5953.099999999999  # M138: -5379.636200000001 -> -5953.099999999999
$ws.Cells.Item(138, 14).Value = -19592.1068  # N138: -20140.0001 -> -19592.1068
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 6162.5  # H110: 3328.3076 -> 6162.5
$ws.Cells.Item(110, 9).Value = 0  # I110: 982.1667 -> 0
$ws.Cells.Item(110, 10).Value = 6162.5  # J110: 5339.2856 -> 6162.5
$ws.Cells.Item(110, 11).Value = 0  # K110: 982.1667 -> 0
$ws.Cells.Item(110, 12).Value = 6162.5  # L110: 5339.2856 -> 6162.5
$ws.Cells.Item(110, 13).ClearContents()  # M110: 1062.8333 -> (removed)
$ws.Cells.Item(110, 14).Value = -10252.5  # N110: -9429.285599999999 -> -10252.5
$ws.Cells.Item(122, 8).Value = 3381.5  # H122: 4914.75 -> 3381.5
$ws.Cells.Item(122, 9).Value = 1871  # I122: 1898.2 -> 1871
$ws.Cells.Item(122, 10).Value = 19997  # J122: 19997.5 -> 19997
$ws.Cells.Item(122, 11).Value = 5613  # K122: 5694.6 -> 5613
$ws.Cells.Item(122, 12).Value = 59991  # L122: 59992.5 -> 59991
$ws.Cells.Item(122, 13).Value = -3163  # M122: -3244.6 -> -3163
$ws.Cells.Item(122, 14).Value = -64891  # N122: -64892.5 -> -64891
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 26002338  # H105: 20002202 -> 26002338
$ws.Cells.Item(105, 9).Value = 1668649.1  # I105: 1113015.9 -> 1668649.1
$ws.Cells.Item(105, 11).Value = 1668649.1  # K105: 1113015.9 -> 1668649.1
$ws.Cells.Item(105, 13).Value = -1666902.1  # M105: -1111268.9 -> -1666902.1
$ws.Cells.Item(129, 8).Value = 199995  # H129: 0 -> 199995
$ws.Cells.Item(129, 10).Value = 199995  # J129: 0 -> 199995
$ws.Cells.Item(129, 12).Value = 199995  # L129: 0 -> 199995
$ws.Cells.Item(129, 14).Value = -209995  # N129: None -> -209995
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 466.2  # H7: 396.66666 -> 466.2
$ws.Cells.Item(7, 9).Value = 115.5  # I7: 200 -> 115.5
$ws.Cells.Item(7, 10).Value = 700  # J7: 436 -> 700
$ws.Cells.Item(7, 11).Value = 115.5  # K7: 200 -> 115.5
$ws.Cells.Item(7, 12).Value = 700  # L7: 436 -> 700
$ws.Cells.Item(7, 13).Value = -2.5  # M7: -87 -> -2.5
$ws.Cells.Item(7, 14).Value = -926  # N7: -662 -> -926
$ws.Cells.Item(16, 8).Value = 2173.4119  # H16: 2216.1177 -> 2173.4119
$ws.Cells.Item(16, 10).Value = 2034.8334  # J16: 2155.8333 -> 2034.8334
$ws.Cells.Item(16, 12).Value = 2034.8334  # L16: 2155.8333 -> 2034.8334
$ws.Cells.Item(16, 14).Value = -2608.8334  # N16: -2729.8333 -> -2608.8334
$ws.Cells.Item(31, 8).Value = 4230.0605  # H31: 4191.6875 -> 4230.0605
$ws.Cells.Item(31, 9).Value = 2825.3914  # I31: 2881.4348 -> 2825.3914
$ws.Cells.Item(31, 10).Value = 7460.8  # J31: 7540.1113 -> 7460.8
$ws.Cells.Item(31, 11).Value = 2825.3914  # K31: 2881.4348 -> 2825.3914
$ws.Cells.Item(31, 12).Value = 7460.8  # L31: 7540.1113 -> 7460.8
$ws.Cells.Item(31, 13).Value = -2530.3914  # M31: -2586.4348 -> -2530.3914
$ws.Cells.Item(31, 14).Value = -8050.8  # N31: -8130.1113 -> -8050.8
$ws.Cells.Item(34, 8).Value = 4230.0605  # H34: 4191.6875 -> 4230.0605
$ws.Cells.Item(34, 9).Value = 2825.3914  # I34: 2881.4348 -> 2825.3914
$ws.Cells.Item(34, 10).Value = 7460.8  # J34: 7540.1113 -> 7460.8
$ws.Cells.Item(34, 11).Value = 2825.3914  # K34: 2881.4348 -> 2825.3914
$ws.Cells.Item(34, 12).Value = 7460.8  # L34: 7540.1113 -> 7460.8
$ws.Cells.Item(34, 13).Value = -2623.3914  # M34: -2679.4348 -> -2623.3914
$ws.Cells.Item(34, 14).Value = -7864.8  # N34: -7944.1113 -> -7864.8
$ws.Cells.Item(58, 8).Value = 2086.2  # H58: 2150.0417 -> 2086.2
$ws.Cells.Item(58, 9).Value = 1078.5  # I58: 1091.1 -> 1078.5
$ws.Cells.Item(58, 10).Value = 2758  # J58: 2906.4285 -> 2758
$ws.Cells.Item(58, 11).Value = 1078.5  # K58: 1091.1 -> 1078.5
$ws.Cells.Item(58, 12).Value = 2758  # L58: 2906.4285 -> 2758
$ws.Cells.Item(58, 13).Value = -875.5  # M58: -888.0999999999999 -> -875.5
$ws.Cells.Item(58, 14).Value = -3164  # N58: -3312.4285 -> -3164
$ws.Cells.Item(103, 8).Value = 4416.7144  # H103: 4677.2856 -> 4416.7144
$ws.Cells.Item(103, 9).Value = 4416.7144  # I103: 4677.2856 -> 4416.7144
$ws.Cells.Item(103, 11).Value = 4416.7144  # K103: 4677.2856 -> 4416.7144
$ws.Cells.Item(103, 13).Value = -3244.7144  # M103: -3505.2856 -> -3244.7144
$ws.Cells.Item(113, 8).Value = 2173.4119  # H113: 2216.1177 -> 2173.4119
$ws.Cells.Item(113, 10).Value = 2034.8334  # J113: 2155.8333 -> 2034.8334
$ws.Cells.Item(113, 12).Value = 2034.8334  # L113: 2155.8333 -> 2034.8334
$ws.Cells.Item(113, 14).Value = -6374.8334  # N113: -6495.8333 -> -6374.8334
$ws.Cells.Item(122, 8).Value = 4654.4614  # H122: 4776.64 -> 4654.4614
$ws.Cells.Item(122, 9).Value = 3928  # I122: 4122 -> 3928
$ws.Cells.Item(122, 11).Value = 11784  # K122: 12366 -> 11784
$ws.Cells.Item(122, 13).Value = -9334  # M122: -9916 -> -9334
$ws.Cells.Item(134, 8).Value = 2311.2974  # H134: 1897.6 -> 2311.2974
$ws.Cells.Item(134, 9).Value = 2047.4375  # I134: 1664 -> 2047.4375
$ws.Cells.Item(134, 11).Value = 6142.3125  # K134: 4992 -> 6142.3125
$ws.Cells.Item(134, 13).Value = -3607.3125  # M134: -2457 -> -3607.3125
$ws.Cells.Item(136, 8).Value = 2086.2  # H136: 2150.0417 -> 2086.2
$ws.Cells.Item(136, 9).Value = 1078.5  # I136: 1091.1 -> 1078.5
$ws.Cells.Item(136, 10).Value = 2758  # J136: 2906.4285 -> 2758
$ws.Cells.Item(136, 11).Value = 3235.5  # K136: 3273.3 -> 3235.5
$ws.Cells.Item(136, 12).Value = 8274  # L136: 8719.2855 -> 8274
$ws.Cells.Item(136, 13).Value = -685.5  # M136: -723.2999999999997 -> -685.5
$ws.Cells.Item(136, 14).Value = -13374  # N136: -13819.2855 -> -13374
$ws.Cells.Item(138, 8).Value = 69999  # H138: 60000 -> 69999
$ws.Cells.Item(138, 10).Value = 69999  # J138: 60000 -> 69999
$ws.Cells.Item(138, 12).Value = 69999  # L138: 60000 -> 69999
$ws.Cells.Item(138, 14).Value = -80279  # N138: -70280 -> -80279
$ws.Cells.Item(141, 8).Value = 331596  # H141: 484271 -> 331596
$ws.Cells.Item(141, 10).Value = 331596  # J141: 484271 -> 331596
$ws.Cells.Item(141, 12).Value = 331596  # L141: 484271 -> 331596
$ws.Cells.Item(141, 14).Value = -341956  # N141: -494631 -> -341956
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 500  # H51: 350 -> 500
$ws.Cells.Item(51, 9).Value = 500  # I51: 350 -> 500
$ws.Cells.Item(51, 11).Value = 1500  # K51: 1050 -> 1500
$ws.Cells.Item(51, 13).Value = -1040  # M51: -590 -> -1040
$ws.Cells.Item(131, 8).Value = 4603.702  # H131: 4459.34 -> 4603.702
$ws.Cells.Item(131, 9).Value = 15570.7  # I131: 15575.6 -> 15570.7
$ws.Cells.Item(131, 10).Value = 1639.6487  # J131: 1680.275 -> 1639.6487
$ws.Cells.Item(131, 11).Value = 46712.10000000001  # K131: 46726.8 -> 46712.10000000001
$ws.Cells.Item(131, 12).Value = 4918.9461  # L131: 5040.825000000001 -> 4918.9461
$ws.Cells.Item(131, 13).Value = -41672.10000000001  # M131: -41686.8 -> -41672.10000000001
$ws.Cells.Item(131, 14).Value = -14998.9461  # N131: -15120.825 -> -14998.9461
$ws.Cells.Item(137, 8).Value = 2549.8948  # H137: 2769.2942 -> 2549.8948
$ws.Cells.Item(137, 9).Value = 2530.5833  # I137: 2728.818 -> 2530.5833
$ws.Cells.Item(137, 10).Value = 2583  # J137: 2843.5 -> 2583
$ws.Cells.Item(137, 11).Value = 7591.749899999999  # K137: 8186.454000000001 -> 7591.749899999999
$ws.Cells.Item(137, 12).Value = 7749  # L137: 8530.5 -> 7749
$ws.Cells.Item(137, 13).Value = -2491.749899999999  # M137: -3086.454000000001 -> -2491.749899999999
$ws.Cells.Item(137, 14).Value = -17949  # N137: -18730.5 -> -17949
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 41667988  # H80: 50001708 -> 41667988
$ws.Cells.Item(80, 9).Value = 55556660  # I80: 71430130 -> 55556660
$ws.Cells.Item(80, 10).Value = 1975.3334  # J80: 2051 -> 1975.3334
$ws.Cells.Item(80, 11).Value = 55556660  # K80: 71430130 -> 55556660
$ws.Cells.Item(80, 12).Value = 1975.3334  # L80: 2051 -> 1975.3334
$ws.Cells.Item(80, 13).Value = -55555662  # M80: -71429132 -> -55555662
$ws.Cells.Item(80, 14).Value = -3971.3334  # N80: -4047 -> -3971.3334
$ws.Cells.Item(83, 8).Value = 41667988  # H83: 50001708 -> 41667988
$ws.Cells.Item(83, 9).Value = 55556660  # I83: 71430130 -> 55556660
$ws.Cells.Item(83, 10).Value = 1975.3334  # J83: 2051 -> 1975.3334
$ws.Cells.Item(83, 11).Value = 277783300  # K83: 357150650 -> 277783300
$ws.Cells.Item(83, 12).Value = 9876.666999999999  # L83: 10255 -> 9876.666999999999
$ws.Cells.Item(83, 13).Value = -277778308  # M83: -357145658 -> -277778308
$ws.Cells.Item(83, 14).Value = -19860.667  # N83: -20239 -> -19860.667
$ws.Cells.Item(122, 8).Value = 4031.2903  # H122: 3850.6667 -> 4031.2903
$ws.Cells.Item(122, 9).Value = 2875.0715  # I122: 2763.5334 -> 2875.0715
$ws.Cells.Item(122, 10).Value = 4983.4707  # J122: 4756.6113 -> 4983.4707
$ws.Cells.Item(122, 11).Value = 8625.2145  # K122: 8290.600199999999 -> 8625.2145
$ws.Cells.Item(122, 12).Value = 14950.4121  # L122: 14269.8339 -> 14950.4121
$ws.Cells.Item(122, 13).Value = -6175.2145  # M122: -5840.600199999999 -> -6175.2145
$ws.Cells.Item(122, 14).Value = -19850.4121  # N122: -19169.8339 -> -19850.4121
$ws.Cells.Item(132, 8).Value = 2371.5557  # H132: 2228.3 -> 2371.5557
$ws.Cells.Item(132, 9).Value = 2268.8  # I132: 1897.5714 -> 2268.8
$ws.Cells.Item(132, 10).Value = 2500  # J132: 3000 -> 2500
$ws.Cells.Item(132, 11).Value = 6806.400000000001  # K132: 5692.7142 -> 6806.400000000001
$ws.Cells.Item(132, 12).Value = 7500  # L132: 9000 -> 7500
$ws.Cells.Item(132, 13).Value = -4276.400000000001  # M132: -3162.7142 -> -4276.400000000001
$ws.Cells.Item(132, 14).Value = -12560  # N132: -14060 -> -12560
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 686.5  # H93: 658 -> 686.5
$ws.Cells.Item(93, 9).Value = 686.5  # I93: 658 -> 686.5
$ws.Cells.Item(93, 11).Value = 686.5  # K93: 658 -> 686.5
$ws.Cells.Item(93, 13).Value = 561.5  # M93: 590 -> 561.5
$ws.Cells.Item(137, 8).Value = 57776.89  # H137: 57940.41 -> 57776.89
$ws.Cells.Item(137, 10).Value = 57776.89  # J137: 57940.41 -> 57776.89
$ws.Cells.Item(137, 12).Value = 57776.89  # L137: 57940.41 -> 57776.89
$ws.Cells.Item(137, 14).Value = -67976.89  # N137: -68140.41 -> -67976.89
$ws.Cells.Item(139, 8).Value = 69999  # H139: 70064.14 -> 69999
$ws.Cells.Item(139, 10).Value = 69999  # J139: 70064.14 -> 69999
$ws.Cells.Item(139, 12).Value = 69999  # L139: 70064.14 -> 69999
$ws.Cells.Item(139, 14).Value = -80279  # N139: -80344.14 -> -80279
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4666.6665  # H81: 4249.75 -> 4666.6665
$ws.Cells.Item(81, 10).Value = 5000  # J81: 3999.5 -> 5000
$ws.Cells.Item(81, 12).Value = 10000  # L81: 7999 -> 10000
$ws.Cells.Item(81, 14).Value = -12122  # N81: -10121 -> -12122
$ws.Cells.Item(84, 8).Value = 4666.6665  # H84: 4249.75 -> 4666.6665
$ws.Cells.Item(84, 10).Value = 5000  # J84: 3999.5 -> 5000
$ws.Cells.Item(84, 12).Value = 50000  # L84: 39995 -> 50000
$ws.Cells.Item(84, 14).Value = -60608  # N84: -50603 -> -60608
$ws.Cells.Item(122, 8).Value = 7813693  # H122: 7576934 -> 7813693
$ws.Cells.Item(122, 9).Value = 892.5  # I122: 881.7826 -> 892.5
$ws.Cells.Item(122, 11).Value = 2677.5  # K122: 2645.3478 -> 2677.5
$ws.Cells.Item(122, 13).Value = -227.5  # M122: -195.3478 -> -227.5
$ws.Cells.Item(126, 8).Value = 2848.2  # H126: 3043.4546 -> 2848.2
$ws.Cells.Item(126, 9).Value = 1797.4286  # I126: 1882.7142 -> 1797.4286
$ws.Cells.Item(126, 10).Value = 5300  # J126: 5074.75 -> 5300
$ws.Cells.Item(126, 11).Value = 5392.2858  # K126: 5648.142599999999 -> 5392.2858
$ws.Cells.Item(126, 12).Value = 15900  # L126: 15224.25 -> 15900
$ws.Cells.Item(126, 13).Value = -2922.2858  # M126: -3178.142599999999 -> -2922.2858
$ws.Cells.Item(126, 14).Value = -20840  # N126: -20164.25 -> -20840
$ws.Cells.Item(139, 8).Value = 80832.336  # H139: 80908.17999999999 -> 80832.336
$ws.Cells.Item(139, 10).Value = 79998.91  # J139: 79999 -> 79998.91
$ws.Cells.Item(139, 12).Value = 79998.91  # L139: 79999 -> 79998.91
$ws.Cells.Item(139, 14).Value = -90278.91  # N139: -90279 -> -90278.91
$ws.Cells.Item(140, 8).Value = 100652.336  # H140: 100652.89 -> 100652.336
$ws.Cells.Item(140, 10).Value = 96935.125  # J140: 96935.75 -> 96935.125
$ws.Cells.Item(140, 12).Value = 96935.125  # L140: 96935.75 -> 96935.125
$ws.Cells.Item(140, 14).Value = -107295.125  # N140: -107295.75 -> -107295.125
